$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row 15 ("Herramienta") holds the per-instrument-column codes. ---
# Cell 4 (criterion column I3): "LC" -> "RUB"
$cLC = $t.Rows.Item(15).Cells.Item(4)
$rLC = $cLC.Range
$rLC.End = $rLC.End - 1           # exclude the end-of-cell mark
$rLC.Text = "RUB"

# Re-fetch the table/row after the mutation above before touching the
# next cell, since cached Cell/Range handles can go stale across edits.
$t = $d.Tables.Item(1)

# Cell 7 (criterion columns I10-I18, merged via gridSpan): "RUB" -> "RUB / EV"
$cRUB = $t.Rows.Item(15).Cells.Item(7)
$rRUB = $cRUB.Range
$rRUB.End = $rRUB.End - 1         # exclude the end-of-cell mark
$rRUB.Text = "RUB / EV"

# Re-fetch again before editing the legend row.
$t = $d.Tables.Item(1)

# Row 17 is the legend explaining each instrument abbreviation; remove the
# "LC -  Lista de control; " entry now that "LC" is no longer used above.
$legendCell = $t.Rows.Item(17).Cells.Item(1)
$legendCell.Range.Find.Execute("LC -  Lista de control; ", $true, $false, $false, $false, $false, $true, 0, $false, "", 1) | Out-Null
